$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MEDICINE")

# Update the last data row (row 11) with the new record details
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "phuc"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "asd"
$ws.Range("E11").Value = "ok"
$ws.Range("F11").Value = 44752

# Re-apply the date format to the ExpireDate column for all data rows
$ws.Range("F2:F11").NumberFormat = "dd/mm/yyyy"
